$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue 'D2' '63.554.78'
Set-TextValue 'E2' '  +1.47%  '
Set-TextValue 'D3' '3.176.87'
Set-TextValue 'E3' '  -0.77%  '
Set-TextValue 'D5' '592.98'
Set-TextValue 'E5' '  -0.44%  '
Set-TextValue 'D6' '135.76'
Set-TextValue 'E6' '  -0.29%  '
Set-TextValue 'E7' '  +0.04%  '
Set-TextValue 'D8' '3.174.03'
Set-TextValue 'E8' '  -0.84%  '
Set-TextValue 'E9' '  +1.66%  '
Set-TextValue 'E10' '  -1.13%  '
Set-TextValue 'D11' '5.36'
Set-TextValue 'E11' '  -0.17%  '
Set-TextValue 'D12' '0.457'
Set-TextValue 'E12' '  +0.30%  '
Set-TextValue 'E13' '  +0.11%  '
Set-TextValue 'D14' '34.68'
Set-TextValue 'E14' '  +3.13%  '
Set-TextValue 'D15' '3.696.27'
Set-TextValue 'E15' '  -0.84%  '
Set-TextValue 'E16' '  -0.46%  '
Set-TextValue 'D17' '3.174.75'
Set-TextValue 'E17' '  -0.66%  '
Set-TextValue 'D18' '63.555.06'
Set-TextValue 'E18' '  +1.28%  '
Set-TextValue 'D19' '6.56'
Set-TextValue 'E19' '  -2.34%  '
Set-TextValue 'D20' '462.62'
Set-TextValue 'E20' '  -0.28%  '
Set-TextValue 'D21' '13.99'
Set-TextValue 'E21' '  +0.00%  '
Set-TextValue 'D22' '0.698'
Set-TextValue 'E22' '  -2.45%  '
Set-TextValue 'D23' '7.70'
Set-TextValue 'E23' '  -0.09%  '
Set-TextValue 'D24' '13.26'
Set-TextValue 'E24' '  -2.31%  '
Set-TextValue 'D25' '83.17'
Set-TextValue 'E25' '  -0.51%  '
Set-TextValue 'E26' '  +0.02%  '
Set-TextValue 'E27' '  -1.65%  '
Set-TextValue 'D28' '0.999'
Set-TextValue 'E28' '  -0.05%  '
Set-TextValue 'E29' '  -1.47%  '
Set-TextValue 'B30' 'ImmutableX'
Set-TextValue 'C30' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D30' '2.07'
Set-TextValue 'E30' '  -0.68%  '
Set-TextValue 'B31' 'RenderToken'
Set-TextValue 'C31' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D31' '7.75'
Set-TextValue 'E31' '  -2.62%  '
Set-TextValue 'D32' '27.36'
Set-TextValue 'E32' '  -0.89%  '
Set-TextValue 'E33' '  -1.52%  '
Set-TextValue 'E34' '  -1.54%  '
Set-TextValue 'E35' '  -2.10%  '
Set-TextValue 'E36' '  +0.22%  '
Set-TextValue 'D37' '51.56'
Set-TextValue 'E37' '  +0.03%  '
Set-TextValue 'D38' '0.0₃0734'
Set-TextValue 'E38' '  +4.63%  '
Set-TextValue 'E39' '  -1.03%  '
Set-TextValue 'D40' '8.14'
Set-TextValue 'E40' '  +0.01%  '
Set-TextValue 'E41' '  -2.18%  '
Set-TextValue 'D42' '2.66'
Set-TextValue 'E42' '  -0.10%  '
Set-TextValue 'D43' '394.44'
Set-TextValue 'E43' '  -6.03%  '
Set-TextValue 'D44' '2.793.80'
Set-TextValue 'E44' '  -7.40%  '
Set-TextValue 'E45' '  -1.18%  '
Set-TextValue 'B46' 'Arweave'
Set-TextValue 'C46' 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
Set-TextValue 'D46' '36.01'
Set-TextValue 'E46' '  -0.54%  '
Set-TextValue 'B47' 'Monero'
Set-TextValue 'C47' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D47' '127.60'
Set-TextValue 'E47' '  +1.78%  '
Set-TextValue 'E48' '  -0.01%  '
Set-TextValue 'E49' '  -2.10%  '
Set-TextValue 'D50' '25.30'
Set-TextValue 'E50' '  -2.84%  '
Set-TextValue 'E51' '  -0.81%  '
